$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a value to a cell while forcing Text storage (matches the
# original inline-string cell type) so numeric-looking strings keep their
# exact formatting (e.g. trailing zeros) instead of being coerced to a Double.
function Set-TextValue([string]$addr, [string]$value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "289.50"
Set-TextValue "D3" "21.26"
Set-TextValue "D4" "6.468"
Set-TextValue "D5" "0.06381"
Set-TextValue "B6" "GateToken"
Set-TextValue "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "3.614"
Set-TextValue "E6" "5GateTokenGT"
Set-TextValue "B7" "FTXToken"
Set-TextValue "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.580"
Set-TextValue "E7" "6FTXTokenFTT"
Set-TextValue "D8" "6.620"
Set-TextValue "D9" "0.8296"
Set-TextValue "D10" "0.01425"
Set-TextValue "D11" "0.1690"
Set-TextValue "D12" "0.08730"
Set-TextValue "D13" "0.03654"
Set-TextValue "D14" "0.03217"
Set-TextValue "D15" "0.09194"
Set-TextValue "D16" "3.717"
Set-TextValue "D17" "0.001657"
Set-TextValue "D18" "0.04747"
Set-TextValue "D19" "0.006117"
Set-TextValue "D20" "0.006302"
Set-TextValue "D21" "0.001073"
Set-TextValue "D22" "0.0001604"
Set-TextValue "D23" "3.765"
Set-TextValue "D25" "0.3360"
Set-TextValue "D26" "0.1263"
Set-TextValue "D28" "0.0002713"
Set-TextValue "D40" "0.04859"
Set-TextValue "D41" "0.007190"
Set-TextValue "D42" "0.004512"
Set-TextValue "D43" "0.1118"
Set-TextValue "D44" "0.01151"
Set-TextValue "D45" "0.00006960"
Set-TextValue "D46" "0.00000000752"
Set-TextValue "D47" "0.8023"
Set-TextValue "D48" "0.007344"
Set-TextValue "D49" "0.00001905"
Set-TextValue "D50" "0.01243"
